# Industrial Controls Notes.xlsx - updated tla/fla and added placeholder
# documents for SEL and DNP protocol cheat sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HW Q's")

# Reword the structured-text question (B35) to call out IEC-61131 explicitly.
$ws.Range("B35").Value = "*Given test statement with ':=' instead of '='* What is wrong with this IEC-61131 structured text statement?"

# Fill in the previously-blank placeholder questions (D8:D17) covering the
# new DNP / Modbus / Telnet protocol cheat-sheet material.
$ws.Range("D8").Value = "Will  all data be collected in a class 1 poll in DNP?"
$ws.Range("D9").Value = "What event data is typically assocaited with a class 0 poll? Class 1? Class 2? Class 3?"
$ws.Range("D10").Value = 'Suggest two reasons that analog values reported by DNP is evaluated as "0" when it is actually measuring a small value.'
$ws.Range("D11").Value = "Are Telnet communications encrypted?"
$ws.Range("D12").Value = "Can Modbus RTU be used over Ethernet?"
$ws.Range("D13").Value = "Is Modbus TCP equivalent to Modbus RTU over Ethernet?"
$ws.Range("D14").Value = "Can Modbus input registers be written to from a client/master?"
$ws.Range("D15").Value = "Which terms are essentially equivalent: (Master/Slave/Client/Server)?"
$ws.Range("D16").Value = "In most vendor implementations, are Modbus maps configurable?"
$ws.Range("D17").Value = "In most vendor implementations, are DNP maps configurable?"

# Reflect the user's last on-screen position: scrolled down with D18 selected.
$ws.Activate()
$ws.Range("D18").Select()
$excel.ActiveWindow.ScrollRow = 7
